$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-10-08 Wednesday"; new = "2025-10-09 Thursday"},
    @{old = "37×51="; new = "41×47="},
    @{old = "17×44="; new = "51×95="},
    @{old = "25×13="; new = "78×19="},
    @{old = "35×21="; new = "62×69="},
    @{old = "71×59="; new = "18×12="},
    @{old = "89×49="; new = "30×91="},
    @{old = "92×24="; new = "39×71="},
    @{old = "19×45="; new = "73×99="},
    @{old = "42×86="; new = "48×74="},
    @{old = "19×53="; new = "25×44="},
    @{old = "42×88="; new = "14×94="},
    @{old = "65×71="; new = "71×82="},
    @{old = "15×75="; new = "63×64="},
    @{old = "66×22="; new = "69×57="},
    @{old = "83×38="; new = "19×94="},
    @{old = "68×74="; new = "90×13="},
    @{old = "65×96="; new = "79×33="},
    @{old = "36×64="; new = "71×41="},
    @{old = "72×28="; new = "15×81="},
    @{old = "86×11="; new = "42×98="},
    @{old = "76×76="; new = "85×17="},
    @{old = "17×31="; new = "44×67="},
    @{old = "28×51="; new = "89×66="},
    @{old = "88×57="; new = "51×68="},
    @{old = "16×73="; new = "64×72="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
